$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Set column H (Industries) values to 0 for rows 20 through 105
$ws.Range("H20:H105").Value = 0
